# Weekly update: insert 3 new price rows (week of 2023-11-09, serial 45239)
# for "Terminal La Palmera de La Serena" / Mandarina / Murcott, ahead of the
# existing 2023-09-04 (45173) rows, shifting everything else down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 1015 (old rows 1015-1047 shift to 1018-1050).
$ws.Range("A1015:A1017").EntireRow.Insert()

# New row 1015: Murcott / Especial
$ws.Range("A1015").Value = 8
$ws.Range("B1015").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1015").Value = "Coquimbo"
$ws.Range("D1015").Value = 45239
$ws.Range("E1015").Value = 4
$ws.Range("F1015").Value = "Fruta"
$ws.Range("G1015").Value = 100102
$ws.Range("H1015").Value = "Cítricos"
$ws.Range("I1015").Value = 100102004
$ws.Range("J1015").Value = "Mandarina"
$ws.Range("K1015").Value = "Murcott"
$ws.Range("L1015").Value = "Especial"
$ws.Range("M1015").Value = 700
$ws.Range("N1015").Value = 7000
$ws.Range("O1015").Value = 7500
$ws.Range("P1015").Value = 7250
$ws.Range("Q1015").Value = "$/bandeja 10 kilos"
$ws.Range("R1015").Value = "Provincia de Limarí"
$ws.Range("S1015").Value = 725
$ws.Range("T1015").Value = 10

# New row 1016: Murcott / Primera
$ws.Range("A1016").Value = 8
$ws.Range("B1016").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1016").Value = "Coquimbo"
$ws.Range("D1016").Value = 45239
$ws.Range("E1016").Value = 4
$ws.Range("F1016").Value = "Fruta"
$ws.Range("G1016").Value = 100102
$ws.Range("H1016").Value = "Cítricos"
$ws.Range("I1016").Value = 100102004
$ws.Range("J1016").Value = "Mandarina"
$ws.Range("K1016").Value = "Murcott"
$ws.Range("L1016").Value = "Primera"
$ws.Range("M1016").Value = 1100
$ws.Range("N1016").Value = 5000
$ws.Range("O1016").Value = 5500
$ws.Range("P1016").Value = 5250
$ws.Range("Q1016").Value = "$/bandeja 10 kilos"
$ws.Range("R1016").Value = "Provincia de Limarí"
$ws.Range("S1016").Value = 525
$ws.Range("T1016").Value = 10

# New row 1017: Murcott / Segunda
$ws.Range("A1017").Value = 8
$ws.Range("B1017").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1017").Value = "Coquimbo"
$ws.Range("D1017").Value = 45239
$ws.Range("E1017").Value = 4
$ws.Range("F1017").Value = "Fruta"
$ws.Range("G1017").Value = 100102
$ws.Range("H1017").Value = "Cítricos"
$ws.Range("I1017").Value = 100102004
$ws.Range("J1017").Value = "Mandarina"
$ws.Range("K1017").Value = "Murcott"
$ws.Range("L1017").Value = "Segunda"
$ws.Range("M1017").Value = 600
$ws.Range("N1017").Value = 3000
$ws.Range("O1017").Value = 3500
$ws.Range("P1017").Value = 3250
$ws.Range("Q1017").Value = "$/bandeja 10 kilos"
$ws.Range("R1017").Value = "Provincia de Limarí"
$ws.Range("S1017").Value = 325
$ws.Range("T1017").Value = 10
